# Generate Report for Handback
# Updates the localization-status report after a handback from de-de:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - The zh-cn / de-de detail sheets record the newly produced target (.md)
#    and handback (.xlf) files together with the handback timestamp
#  - Column widths are widened to fit the new, longer file-name values

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6a6d28f314d2a987b96cb67ab7c054deb83c82d5/e2e/16f4bed0-8b99-4d3d-ac91-84247a6c6002.md"
$mdName    = "16f4bed0-8b99-4d3d-ac91-84247a6c6002.md"
$zhXlf     = "16f4bed0-8b99-4d3d-ac91-84247a6c6002.b95ab77519e70a783ed97edf66d60250b2027203.zh-cn.xlf"
$deXlf     = "16f4bed0-8b99-4d3d-ac91-84247a6c6002.b95ab77519e70a783ed97edf66d60250b2027203.de-de.xlf"
$status    = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: status column for both languages, wider columns
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $status
$wsOverview.Range("F2").Value = $status
$wsOverview.Columns("E").ColumnWidth = 29.166666666666664
$wsOverview.Columns("F").ColumnWidth = 29.166666666666664

# ---------------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $status

# Newly generated target file (hyperlinked, same as the source-file link style)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $targetUrl, "", "", $mdName)
$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("I2").Font.Color = 15570276

# Newly generated handback xliff file
$wsZh.Range("J2").Value = $zhXlf

# Latest Handback DateTime (existing "0001-01-01 00:00:00" placeholder is
# replaced in place with the real handback timestamp)
$wsZh.Range("K2").Value = "2016-09-05 23:12:44"

$wsZh.Columns("C").ColumnWidth = 29.166666666666664
$wsZh.Columns("I").ColumnWidth = 39.16666666666667
$wsZh.Columns("J").ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $status

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $targetUrl, "", "", $mdName)
$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("I2").Font.Color = 15570276

$wsDe.Range("J2").Value = $deXlf

# de-de is the language that was actually handed back, so it gets the newly
# minted handback timestamp
$wsDe.Range("K2").Value = "2016-09-05 23:12:51"

$wsDe.Columns("C").ColumnWidth = 29.166666666666664
$wsDe.Columns("I").ColumnWidth = 39.16666666666667
$wsDe.Columns("J").ColumnWidth = 39.16666666666667
